# Final Uploading of re-labelled parcel and parcellation data.
# Adds a new "Goochland" summary sheet (re-labelled codes) positioned
# right after "Powhatan" and before "Goochland 2018".

$wb = $excel.ActiveWorkbook

# Insert the new worksheet directly after the "Powhatan" tab so it lands
# between "Powhatan" and "Goochland 2018", matching the target layout.
$afterSheet = $wb.Worksheets.Item("Powhatan")
$goochland = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$goochland.Name = "Goochland"

# Header row
$goochland.Range("B2").Value = "gooch_value"
$goochland.Range("C2").Value = "gooch_code"

# Re-labelled value/code pairs
$goochland.Range("B3").Value = "other"
$goochland.Range("C3").Value = 7

$goochland.Range("B4").Value = "comm/ind"
$goochland.Range("C4").Value = 4

$goochland.Range("B5").Value = "res sub"
$goochland.Range("C5").Value = 2

$goochland.Range("B6").Value = "ag20-100"
$goochland.Range("C6").Value = 5

$goochland.Range("B7").Value = "ag100+"
$goochland.Range("C7").Value = 6

$goochland.Range("B8").Value = "res urb"
$goochland.Range("C8").Value = 1

$goochland.Range("B9").Value = "multi"
$goochland.Range("C9").Value = 3

$goochland.Range("B10").Value = "NA"
$goochland.Range("C10").Value = 0

# Column widths to roughly match the authored layout
$goochland.Columns.Item(2).ColumnWidth = 21.75
$goochland.Columns.Item(3).ColumnWidth = 21.25

# Leave the selection where the author last left it on this sheet
$goochland.Range("B15").Select() | Out-Null
